$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GA")

# Update C2 value
$ws.Range("C2").Value = 5952

# Add new cells E2, F2, G2
$ws.Range("E2").Value = 1001
$ws.Range("F2").Value = 4951
$ws.Range("G2").Formula = "=E2+F2"

# Apply wrap text style to F2 (new cellXf with wrapText=true)
$ws.Range("F2").WrapText = $true

# Row height grows to fit wrapped text
$ws.Rows.Item(2).RowHeight = 14.9

# Select C2 on the GA sheet
$ws.Range("C2").Select()
